$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New TPM-derived values per row/column, per commit 'update scripts wuth new tpm'
$updates = @{
    2 = @{ "E"=2; "F"=0.6666666666666666; "G"=0.04422566666666666; "H"=0.132677; "I"=0.00399376563542826; "J"=0.003993765635428259; "M"=123.2806423333333; "N"=369.841927; "O"=0.6241574062367528; "P"=0.6241574062367526; "Q"=5.452168594286555; "R"=49.069517348579; "S"=0.002492738400126379; "T"=0.002492738400126379 }
    3 = @{ "E"=2; "F"=0.6666666666666666; "G"=0.04422566666666666; "H"=0.132677; "I"=0.00399376563542826; "J"=0.003993765635428259; "O"=0.2392728888301323; "P"=0.2392728888301322; "Q"=2.090107586497222; "R"=18.81096827847499; "S"=0.0009555998408994287; "T"=0.0009555998408994283 }
    4 = @{ "E"=2; "F"=0.6666666666666666; "G"=0.04422566666666666; "H"=0.132677; "I"=0.00399376563542826; "J"=0.003993765635428259; "O"=0.136569704933115; "P"=0.136569704933115; "Q"=1.192969992388222; "R"=10.736729931494; "S"=0.0005454273944024522; "T"=0.000545427394402452 }
    5 = @{ "G"=8.111092999999999; "I"=0.732466165706853; "J"=0.7324661657068529; "M"=123.2806423333333; "N"=369.841927; "O"=0.6241574062367528; "P"=0.6241574062367526; "Q"=999.9407550654034; "R"=8999.466795588633; "S"=0.4571741821437689; "T"=0.4571741821437688 }
    6 = @{ "G"=8.111092999999999; "I"=0.732466165706853; "J"=0.7324661657068529; "O"=0.2392728888301323; "P"=0.2392728888301322; "Q"=383.3307283270915; "R"=3449.976554943824; "S"=0.1752592954390091; "T"=0.175259295439009 }
    7 = @{ "G"=8.111092999999999; "I"=0.732466165706853; "J"=0.7324661657068529; "O"=0.136569704933115; "P"=0.136569704933115; "Q"=218.7935487191486; "S"=0.1000326881240751; "T"=0.100032688124075 }
    8 = @{ "H"=8.755071999999998; "I"=0.2635400686577188; "J"=0.2635400686577188; "M"=123.2806423333333; "N"=369.841927; "O"=0.6241574062367528; "P"=0.6241574062367526; "Q"=359.776966611527; "R"=3237.992699503744; "S"=0.1644904856928575; "T"=0.1644904856928575 }
    9 = @{ "H"=8.755071999999998; "I"=0.2635400686577188; "J"=0.2635400686577188; "O"=0.2392728888301323; "P"=0.2392728888301322; "S"=0.06305799355022378; "T"=0.06305799355022378 }
    10 = @{ "H"=8.755071999999998; "I"=0.2635400686577188; "J"=0.2635400686577188; "O"=0.136569704933115; "P"=0.136569704933115; "Q"=78.72154312502043; "R"=708.4938881251838; "S"=0.03599158941463754; "T"=0.03599158941463754 }
}

foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row].Keys) {
        $ws.Range("$col$row").Value = $updates[$row][$col]
    }
}
